$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = "La Thị Hồng Nhung"
$ws.Range("C4").Value = "B20DCDT158"
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 7.3
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 7
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 6.75
$ws.Range("O4").Value = 6.428571428571429
$ws.Range("P4").Value = 6.285714285714286
$ws.Range("Q4").Value = 6.571428571428571
$ws.Range("R4").Value = 6.839285714285714
$ws.Range("S4").Value = 8
$ws.Range("T4").Value = 6.571428571428571
$ws.Range("U4").Value = 6.428571428571429
$ws.Range("V4").Value = 5.8
$ws.Range("W4").Value = 6.571428571428571
$ws.Range("X4").Value = 6.674285714285714
$ws.Range("Y4").Value = 7.059749999999999

# Row 5
$ws.Range("B5").Value = "Nguyễn Trung Tuấn"
$ws.Range("C5").Value = "B20DCDT190"
$ws.Range("N5").Value = 3.625
$ws.Range("O5").Value = 4.142857142857143
$ws.Range("P5").Value = 4.142857142857143
$ws.Range("Q5").Value = 4.142857142857143
$ws.Range("R5").Value = 2.675595238095239
$ws.Range("T5").Value = 4.142857142857143
$ws.Range("U5").Value = 4.142857142857143
$ws.Range("V5").Value = 5.8
$ws.Range("W5").Value = 4.142857142857143
$ws.Range("X5").Value = 3.645714285714285
$ws.Range("Y5").Value = 2.212458333333333

# Row 6
$ws.Range("B6").Value = "Phạm Thế Anh"
$ws.Range("C6").Value = "B20DCDT017"
$ws.Range("D6").Value = "D20DTMT1"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5
$ws.Range("N6").Value = 3.375
$ws.Range("O6").Value = 3.857142857142857
$ws.Range("P6").Value = 3.714285714285714
$ws.Range("Q6").Value = 3.857142857142857
$ws.Range("R6").Value = 2.717261904761905
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 3.857142857142857
$ws.Range("U6").Value = 3.857142857142857
$ws.Range("V6").Value = 5
$ws.Range("W6").Value = 3.857142857142857
$ws.Range("X6").Value = 3.514285714285715
$ws.Range("Y6").Value = 2.281041666666667

# Row 7
$ws.Range("B7").Value = "Nguyễn Tiến Duy"
$ws.Range("C7").Value = "B20DCDT037"
$ws.Range("N7").Value = 3.125
$ws.Range("O7").Value = 3.571428571428572
$ws.Range("P7").Value = 3.571428571428572
$ws.Range("Q7").Value = 3.571428571428572
$ws.Range("R7").Value = 2.306547619047619
$ws.Range("T7").Value = 3.571428571428572
$ws.Range("U7").Value = 3.571428571428572
$ws.Range("V7").Value = 5
$ws.Range("W7").Value = 3.571428571428572
$ws.Range("X7").Value = 3.142857142857143
$ws.Range("Y7").Value = 1.907291666666667

# Row 8
$ws.Range("B8").Value = "Lê Sỹ Sang"
$ws.Range("C8").Value = "B20DCDT175"
$ws.Range("N8").Value = 3.125
$ws.Range("O8").Value = 3.571428571428572
$ws.Range("P8").Value = 3.571428571428572
$ws.Range("Q8").Value = 3.571428571428572
$ws.Range("R8").Value = 2.306547619047619
$ws.Range("T8").Value = 3.571428571428572
$ws.Range("U8").Value = 3.571428571428572
$ws.Range("V8").Value = 5
$ws.Range("W8").Value = 3.571428571428572
$ws.Range("X8").Value = 3.142857142857143
$ws.Range("Y8").Value = 1.907291666666667
